$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 with the combined/serialized text
$ws.Range("A2").Value = "('Servo', ['Token Artifact Creature — Servo', '1/1', 'Thopter', 'Token Artifact Creature — Thopter', 'Flying', '1/1'])"

# Remove now-obsolete rows 3 through 8
$ws.Range("A3:A8").ClearContents()
